$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 315, shifting rows 315:361 down to 316:362
$ws.Rows(315).Insert()

# Populate the newly inserted row 315 with the new weekly price record
$ws.Cells.Item(315, 1).Value = 10
$ws.Cells.Item(315, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(315, 3).Value = "La Araucanía"
$ws.Cells.Item(315, 4).Value = 44504
$ws.Cells.Item(315, 5).Value = 9
$ws.Cells.Item(315, 6).Value = "Fruta"
$ws.Cells.Item(315, 7).Value = 100108
$ws.Cells.Item(315, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(315, 9).Value = 100108006
$ws.Cells.Item(315, 10).Value = "Plátano"
$ws.Cells.Item(315, 11).Value = "Sin especificar"
$ws.Cells.Item(315, 12).Value = "Pintón"
$ws.Cells.Item(315, 13).Value = 1720
$ws.Cells.Item(315, 14).Value = 18000
$ws.Cells.Item(315, 15).Value = 21000
$ws.Cells.Item(315, 16).Value = 19953
$ws.Cells.Item(315, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(315, 18).Value = "Ecuador"
$ws.Cells.Item(315, 19).Value = 998
$ws.Cells.Item(315, 20).Value = 20
